$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 361.42856
$ws.Range("J2").Value = 350
$ws.Range("L2").Value = 350
$ws.Range("N2").Value = -576
$ws.Range("L6").Value = 600
$ws.Range("I6").Value = 221.14285
$ws.Range("K6").Value = 663.4285500000001
$ws.Range("H6").Value = 214.8
$ws.Range("J6").Value = 200
$ws.Range("N6").Value = -824
$ws.Range("M6").Value = -551.4285500000001
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = 0
$ws.Range("I9").Value = 120
$ws.Range("K9").Value = 120
$ws.Range("H9").Value = 120
$ws.Range("M9").Value = 49
$ws.Range("I98").Value = 1670.8889
$ws.Range("K98").Value = 1670.8889
$ws.Range("H98").Value = 1670.8889
$ws.Range("M98").Value = -172.8888999999999
$ws.Range("I122").Value = 1670.8889
$ws.Range("H122").Value = 1670.8889
$ws.Range("K122").Value = 5012.6667
$ws.Range("M122").Value = -2562.6667
$ws.Range("I127").Value = 8408.286
$ws.Range("M127").Value = -20264.858
$ws.Range("K127").Value = 25224.858
$ws.Range("H127").Value = 8408.286
$ws.Range("H138").Value = 2000.7937
$ws.Range("J138").Value = 2094.2144
$ws.Range("N138").Value = -16562.6432
$ws.Range("L138").Value = 6282.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K8").Value = 50002500
$ws.Range("N8").Value = -11787
$ws.Range("M8").Value = -50002356
$ws.Range("H8").Value = 25007000
$ws.Range("L8").Value = 11499
$ws.Range("J8").Value = 11499
$ws.Range("I8").Value = 50002500
$ws.Range("M13").Value = -500856
$ws.Range("I13").Value = 501000
$ws.Range("H13").Value = 336666.34
$ws.Range("K13").Value = 501000
$ws.Range("H30").Value = 498
$ws.Range("N30").ClearContents()
$ws.Range("J30").Value = 0
$ws.Range("I30").Value = 498
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -348
$ws.Range("K30").Value = 498
$ws.Range("H39").Value = 4666.6665
$ws.Range("I39").Value = 4666.6665
$ws.Range("M39").Value = -4146.6665
$ws.Range("K39").Value = 4666.6665
$ws.Range("K45").Value = 1997
$ws.Range("I45").Value = 1997
$ws.Range("H45").Value = 3249.25
$ws.Range("M45").Value = -1620
$ws.Range("H76").Value = 39000
$ws.Range("L76").Value = 39000
$ws.Range("J76").Value = 39000
$ws.Range("N76").Value = -39676
$ws.Range("L79").Value = 39000
$ws.Range("N79").Value = -41340
$ws.Range("J79").Value = 39000
$ws.Range("H79").Value = 39000
$ws.Range("I122").Value = 1258512.8
$ws.Range("H122").Value = 848543.25
$ws.Range("K122").Value = 3775538.4
$ws.Range("M122").Value = -3773088.4
$ws.Range("H132").Value = 4249.5454
$ws.Range("M132").Value = -6468.5
$ws.Range("K132").Value = 8998.5
$ws.Range("I132").Value = 2999.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value = 380
$ws.Range("I10").Value = 10
$ws.Range("J10").Value = 750
$ws.Range("K10").Value = 10
$ws.Range("M10").Value = 130
$ws.Range("N10").Value = -1030
$ws.Range("L10").Value = 750
$ws.Range("J12").Value = 1999.6666
$ws.Range("M12").Value = -229.5
$ws.Range("N12").Value = -2335.6666
$ws.Range("H12").Value = 1084.1428
$ws.Range("K12").Value = 397.5
$ws.Range("L12").Value = 1999.6666
$ws.Range("I12").Value = 397.5
$ws.Range("J47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("L47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("H99").Value = 1877.25
$ws.Range("N99").Value = -5328.6667
$ws.Range("M99").Value = -106
$ws.Range("L99").Value = 2332.6667
$ws.Range("I99").Value = 1604
$ws.Range("K99").Value = 1604
$ws.Range("J99").Value = 2332.6667
$ws.Range("I107").Value = 1745.5
$ws.Range("K107").Value = 1745.5
$ws.Range("H107").Value = 3083.3333
$ws.Range("M107").Value = 174.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I3").Value = 22999.666
$ws.Range("L3").Value = 28327.334
$ws.Range("M3").Value = -22886.666
$ws.Range("K3").Value = 22999.666
$ws.Range("H3").Value = 25663.5
$ws.Range("J3").Value = 28327.334
$ws.Range("N3").Value = -28553.334
$ws.Range("J12").Value = 3492.5
$ws.Range("M12").Value = -380
$ws.Range("N12").Value = -3832.5
$ws.Range("H12").Value = 2511.6667
$ws.Range("K12").Value = 550
$ws.Range("L12").Value = 3492.5
$ws.Range("I12").Value = 550
$ws.Range("H88").Value = 17319.584
$ws.Range("L88").Value = 17319.584
$ws.Range("J88").Value = 17319.584
$ws.Range("N88").Value = -18131.584
$ws.Range("J91").Value = 17319.584
$ws.Range("N91").Value = -20127.584
$ws.Range("H91").Value = 17319.584
$ws.Range("L91").Value = 17319.584
$ws.Range("N92").Value = -18192.333
$ws.Range("J92").Value = 13200.333
$ws.Range("L92").Value = 13200.333
$ws.Range("H92").Value = 13200.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J55").Value = 7168.3335
$ws.Range("N55").Value = -21859.0005
$ws.Range("L55").Value = 21505.0005
$ws.Range("H55").Value = 5401.25
$ws.Range("J75").Value = 766.3333
$ws.Range("H75").Value = 621
$ws.Range("L75").Value = 2298.9999
$ws.Range("N75").Value = -4294.9999
$ws.Range("L78").Value = 6896.9997
$ws.Range("N78").Value = -16880.9997
$ws.Range("H78").Value = 621
$ws.Range("J78").Value = 766.3333
$ws.Range("J131").Value = 1789.9688
$ws.Range("H131").Value = 1781.4412
$ws.Range("N131").Value = -15449.9064
$ws.Range("L131").Value = 5369.9064
$ws.Range("N137").Value = -21299.4999
$ws.Range("J137").Value = 3699.8333
$ws.Range("H137").Value = 4078.0908
$ws.Range("L137").Value = 11099.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N9").Value = -790
$ws.Range("L9").Value = 450
$ws.Range("I9").Value = 371
$ws.Range("K9").Value = 371
$ws.Range("J9").Value = 450
$ws.Range("H9").Value = 390.75
$ws.Range("M9").Value = -201
$ws.Range("H10").Value = 300
$ws.Range("J10").Value = 500
$ws.Range("N10").Value = -838
$ws.Range("L10").Value = 500
$ws.Range("J12").Value = 2500
$ws.Range("N12").Value = -2780
$ws.Range("H12").Value = 3400
$ws.Range("L12").Value = 2500
$ws.Range("M80").Value = -1875.8
$ws.Range("I80").Value = 2873.8
$ws.Range("L80").Value = 3977.8462
$ws.Range("J80").Value = 3977.8462
$ws.Range("H80").Value = 3497.8262
$ws.Range("K80").Value = 2873.8
$ws.Range("N80").Value = -5973.8462
$ws.Range("M83").Value = -9377
$ws.Range("L83").Value = 19889.231
$ws.Range("H83").Value = 3497.8262
$ws.Range("J83").Value = 3977.8462
$ws.Range("N83").Value = -29873.231
$ws.Range("I83").Value = 2873.8
$ws.Range("K83").Value = 14369
$ws.Range("N92").Value = -17992.333
$ws.Range("J92").Value = 14248.333
$ws.Range("L92").Value = 14248.333
$ws.Range("H92").Value = 14248.333
$ws.Range("I122").Value = 2119.111
$ws.Range("H122").Value = 73819.14
$ws.Range("K122").Value = 6357.333
$ws.Range("M122").Value = -3907.333
$ws.Range("H132").Value = 2009.4286
$ws.Range("M132").Value = -2450.6666
$ws.Range("K132").Value = 4980.6666
$ws.Range("I132").Value = 1660.2222

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K4").Value = 5672.3335
$ws.Range("L4").Value = 25000
$ws.Range("J4").Value = 25000
$ws.Range("M4").Value = -5559.3335
$ws.Range("I4").Value = 5672.3335
$ws.Range("N4").Value = -25226
$ws.Range("H4").Value = 10504.25
$ws.Range("M13").Value = -49999960
$ws.Range("I13").Value = 50000100
$ws.Range("H13").Value = 33334400
$ws.Range("K13").Value = 50000100
$ws.Range("M28").Value = -5440.3335
$ws.Range("J28").Value = 25000
$ws.Range("K28").Value = 5672.3335
$ws.Range("H28").Value = 10504.25
$ws.Range("N28").Value = -25464
$ws.Range("I28").Value = 5672.3335
$ws.Range("L28").Value = 25000
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("N37").Value = -25214
$ws.Range("M37").Value = -5565.3335
$ws.Range("I37").Value = 5672.3335
$ws.Range("L37").Value = 25000
$ws.Range("J37").Value = 25000
$ws.Range("H37").Value = 10504.25
$ws.Range("K37").Value = 5672.3335
$ws.Range("I122").Value = 3196.6667
$ws.Range("H122").Value = 3133.75
$ws.Range("K122").Value = 9590.000100000001
$ws.Range("M122").Value = -7140.000100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("I7").Value = 0
$ws.Range("N7").Value = -478.5
$ws.Range("L7").Value = 252.5
$ws.Range("H7").Value = 252.5
$ws.Range("J7").Value = 252.5
$ws.Range("H17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("K17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("N64").Value = -43049
$ws.Range("H64").Value = 42553
$ws.Range("J64").Value = 42553
$ws.Range("L64").Value = 42553
$ws.Range("H67").Value = 42553
$ws.Range("J67").Value = 42553
$ws.Range("N67").Value = -44269
$ws.Range("L67").Value = 42553
$ws.Range("H132").Value = 1900.6364
$ws.Range("M132").Value = -3412.1
$ws.Range("K132").Value = 5942.1
$ws.Range("I132").Value = 1980.7
